# edit.ps1
# Applies the "custom accuracy + 데이터 1000개" commit:
#  1. Replace the 4 data rows (rows 2-5) with a new block of sensor readings
#     (new timestamps + new J1..J33 flow values).
#  2. Delete the now-superfluous 5th data row (old row 6), shrinking the
#     used range from A1:AH6 to A1:AH5.
#  3. Widen several data columns from 7/8 to 8/9 characters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Overwrite rows 2-5 with the new dataset values ---
# Row 2
$ws.Range("A2").Value = 45115.50694444445
$ws.Range("B2").Value = 14.835
$ws.Range("C2").Value = 9.791
$ws.Range("D2").Value = 3.698
$ws.Range("E2").Value = 32.243
$ws.Range("F2").Value = 24.166
$ws.Range("G2").Value = 11.51
$ws.Range("H2").Value = 34.958
$ws.Range("I2").Value = 18.033
$ws.Range("J2").Value = 7.29
$ws.Range("K2").Value = 10.735
$ws.Range("L2").Value = 12.533
$ws.Range("M2").Value = 13.25
$ws.Range("N2").Value = 3.739
$ws.Range("O2").Value = 11.655
$ws.Range("P2").Value = 16.06
$ws.Range("Q2").Value = 10.282
$ws.Range("R2").Value = 3.096
$ws.Range("S2").Value = 1.74
$ws.Range("T2").Value = 170.025
$ws.Range("U2").Value = 32.298
$ws.Range("V2").Value = 10.758
$ws.Range("W2").Value = 20.812
$ws.Range("X2").Value = 10.713
$ws.Range("Y2").Value = 2.837
$ws.Range("Z2").Value = 18.288
$ws.Range("AA2").Value = 9.502000000000001
$ws.Range("AB2").Value = 8.641999999999999
$ws.Range("AC2").Value = 10.303
$ws.Range("AD2").Value = 12.679
$ws.Range("AE2").Value = 3.311
$ws.Range("AF2").Value = 31.418
$ws.Range("AG2").Value = 5.68
$ws.Range("AH2").Value = 13.449

# Row 3
$ws.Range("A3").Value = 45115.51388888889
$ws.Range("B3").Value = 12.438
$ws.Range("C3").Value = 8.742000000000001
$ws.Range("D3").Value = 1.644
$ws.Range("E3").Value = 27.34
$ws.Range("F3").Value = 21.233
$ws.Range("G3").Value = 9.670999999999999
$ws.Range("H3").Value = 37.64
$ws.Range("I3").Value = 15.125
$ws.Range("J3").Value = 6.446
$ws.Range("K3").Value = 9.346
$ws.Range("L3").Value = 10.795
$ws.Range("M3").Value = 11.458
$ws.Range("N3").Value = 3.141
$ws.Range("O3").Value = 9.775
$ws.Range("P3").Value = 13.693
$ws.Range("Q3").Value = 8.586
$ws.Range("R3").Value = 1.431
$ws.Range("S3").Value = 0.907
$ws.Range("T3").Value = 141.465
$ws.Range("U3").Value = 27.369
$ws.Range("V3").Value = 9.023
$ws.Range("W3").Value = 17.954
$ws.Range("X3").Value = 9.353999999999999
$ws.Range("Y3").Value = 1.907
$ws.Range("Z3").Value = 18.502
$ws.Range("AA3").Value = 7.97
$ws.Range("AB3").Value = 7.234
$ws.Range("AC3").Value = 8.537000000000001
$ws.Range("AD3").Value = 11.058
$ws.Range("AE3").Value = 1.246
$ws.Range("AF3").Value = 34.625
$ws.Range("AG3").Value = 4.876
$ws.Range("AH3").Value = 11.281

# Row 4
$ws.Range("A4").Value = 45115.52083333334
$ws.Range("B4").Value = 22.531
$ws.Range("C4").Value = 16.557
$ws.Range("D4").Value = 1.514
$ws.Range("E4").Value = 49.272
$ws.Range("F4").Value = 39.785
$ws.Range("G4").Value = 17.641
$ws.Range("H4").Value = 66.913
$ws.Range("I4").Value = 27.341
$ws.Range("J4").Value = 12.058
$ws.Range("K4").Value = 17.777
$ws.Range("L4").Value = 19.667
$ws.Range("M4").Value = 20.859
$ws.Range("N4").Value = 5.676
$ws.Range("O4").Value = 17.67
$ws.Range("P4").Value = 25.073
$ws.Range("Q4").Value = 15
$ws.Range("R4").Value = 1.056
$ws.Range("S4").Value = 0.993
$ws.Range("T4").Value = 261.643
$ws.Range("U4").Value = 49.381
$ws.Range("V4").Value = 16.31
$ws.Range("W4").Value = 33.055
$ws.Range("X4").Value = 17.331
$ws.Range("Y4").Value = 2.736
$ws.Range("Z4").Value = 32.893
$ws.Range("AA4").Value = 14.407
$ws.Range("AB4").Value = 12.817
$ws.Range("AC4").Value = 15.087
$ws.Range("AD4").Value = 20.471
$ws.Range("AE4").Value = 0.766
$ws.Range("AF4").Value = 60.925
$ws.Range("AG4").Value = 9.111000000000001
$ws.Range("AH4").Value = 20.392

# Row 5
$ws.Range("A5").Value = 45115.52777777778
$ws.Range("B5").Value = 10.04
$ws.Range("C5").Value = 7.29
$ws.Range("D5").Value = 0.89
$ws.Range("E5").Value = 22.07
$ws.Range("F5").Value = 17.52
$ws.Range("G5").Value = 7.83
$ws.Range("H5").Value = 34.57
$ws.Range("I5").Value = 12.22
$ws.Range("J5").Value = 5.34
$ws.Range("K5").Value = 7.76
$ws.Range("L5").Value = 8.779999999999999
$ws.Range("M5").Value = 9.35
$ws.Range("N5").Value = 2.54
$ws.Range("O5").Value = 7.9
$ws.Range("P5").Value = 11.17
$ws.Range("Q5").Value = 6.83
$ws.Range("R5").Value = 0.73
$ws.Range("S5").Value = 0.53
$ws.Range("T5").Value = 112.87
$ws.Range("U5").Value = 22.22
$ws.Range("V5").Value = 7.29
$ws.Range("W5").Value = 14.75
$ws.Range("X5").Value = 7.69
$ws.Range("Y5").Value = 1.35
$ws.Range("Z5").Value = 16.31
$ws.Range("AA5").Value = 6.44
$ws.Range("AB5").Value = 5.79
$ws.Range("AC5").Value = 6.81
$ws.Range("AD5").Value = 9.08
$ws.Range("AE5").Value = 0.55
$ws.Range("AF5").Value = 31.75
$ws.Range("AG5").Value = 4.01
$ws.Range("AH5").Value = 9.109999999999999

# --- 2) Drop the old 6th row; used range becomes A1:AH5 ---
$ws.Rows.Item(6).Delete()

# --- 3) Widen columns (OOXML <col width> = Excel ColumnWidth + 0.83) ---
# target width 8  -> ColumnWidth 7.17
# target width 9  -> ColumnWidth 8.17
$ws.Columns.Item(2).ColumnWidth = 7.17  # B
$ws.Columns.Item(3).ColumnWidth = 7.17  # C
$ws.Columns.Item(6).ColumnWidth = 7.17  # F
$ws.Columns.Item(7).ColumnWidth = 7.17  # G
$ws.Columns.Item(9).ColumnWidth = 7.17  # I
$ws.Columns.Item(10).ColumnWidth = 7.17  # J
$ws.Columns.Item(11).ColumnWidth = 7.17  # K
$ws.Columns.Item(12).ColumnWidth = 7.17  # L
$ws.Columns.Item(13).ColumnWidth = 7.17  # M
$ws.Columns.Item(15).ColumnWidth = 7.17  # O
$ws.Columns.Item(16).ColumnWidth = 7.17  # P
$ws.Columns.Item(17).ColumnWidth = 7.17  # Q
$ws.Columns.Item(20).ColumnWidth = 8.17  # T
$ws.Columns.Item(22).ColumnWidth = 7.17  # V
$ws.Columns.Item(23).ColumnWidth = 7.17  # W
$ws.Columns.Item(24).ColumnWidth = 7.17  # X
$ws.Columns.Item(26).ColumnWidth = 7.17  # Z
$ws.Columns.Item(27).ColumnWidth = 7.17  # AA
$ws.Columns.Item(28).ColumnWidth = 7.17  # AB
$ws.Columns.Item(29).ColumnWidth = 7.17  # AC
$ws.Columns.Item(30).ColumnWidth = 7.17  # AD
$ws.Columns.Item(34).ColumnWidth = 7.17  # AH
